$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect updated "through" date
$ws.Name = "Through 2021-10-07"

# Row 11 (September) updates
$ws.Range("T11").Value = 4
$ws.Range("U11").Value = 174
$ws.Range("V11").Value = 0.0225

# Row 12 (October) updates
$ws.Range("A12").Value = "October (through 10-07)"
$ws.Range("F12").Value = 10
$ws.Range("L12").Value = 16
$ws.Range("R12").Value = 35
$ws.Range("U12").Value = 51

# Row 13 (Total) updates
$ws.Range("F13").Value = 393
$ws.Range("G13").Value = 0.1048
$ws.Range("L13").Value = 503
$ws.Range("M13").Value = 0.1082
$ws.Range("R13").Value = 883
$ws.Range("S13").Value = 0.0566
$ws.Range("T13").Value = 80
$ws.Range("U13").Value = 1220
$ws.Range("V13").Value = 0.0615
